# Added Week 15 simulations
# Update the "R" (road) row (row 3) target-depth stats on both the
# OFF and DEF sheets with the latest simulation results.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 177
$wsOff.Range("C3").Value = 114
$wsOff.Range("D3").Value = 36
$wsOff.Range("E3").Value = 14
$wsOff.Range("F3").Value = 3
$wsOff.Range("G3").Value = 6

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 177
$wsDef.Range("C3").Value = 129
$wsDef.Range("D3").Value = 41
$wsDef.Range("E3").Value = 22
